$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.8247399999999999
$ws.Range("H2").Value = 2.47422
$ws.Range("I2").Value = 0.00290110066097835
$ws.Range("J2").Value = 0.0029096007056883
$ws.Range("M2").Value = 0.8247399999999999
$ws.Range("N2").Value = 2.47422
$ws.Range("O2").Value = 0.00290110066097835
$ws.Range("P2").Value = 0.0029096007056883
$ws.Range("Q2").Value = 0.6801960675999998
$ws.Range("R2").Value = 6.121764608399999
$ws.Range("S2").Value = 0.000008416385045129019
$ws.Range("T2").Value = 0.000008465776266541853

$ws.Range("G3").Value = 0.8247399999999999
$ws.Range("H3").Value = 2.47422
$ws.Range("I3").Value = 0.00290110066097835
$ws.Range("J3").Value = 0.0029096007056883
$ws.Range("O3").Value = 0.001056260848766244
$ws.Range("P3").Value = 0.001059355627434439
$ws.Range("Q3").Value = 0.2476523773733333
$ws.Range("R3").Value = 2.22887139636
$ws.Range("S3").Value = 0.000003064319046521302
$ws.Range("T3").Value = 0.000003082301881158114

$ws.Range("G4").Value = 0.8247399999999999
$ws.Range("H4").Value = 2.47422
$ws.Range("I4").Value = 0.00290110066097835
$ws.Range("J4").Value = 0.0029096007056883
$ws.Range("M4").Value = 120.5783256666667
$ws.Range("N4").Value = 361.734977
$ws.Range("O4").Value = 0.4241456220035762
$ws.Range("P4").Value = 0.4253883423266084
$ws.Range("Q4").Value = 99.44576831032666
$ws.Range("R4").Value = 895.0119147929399
$ws.Range("S4").Value = 0.001230489144345648
$ws.Range("T4").Value = 0.001237710221025076

$ws.Range("G5").Value = 0.8247399999999999
$ws.Range("H5").Value = 2.47422
$ws.Range("I5").Value = 0.00290110066097835
$ws.Range("J5").Value = 0.0029096007056883
$ws.Range("M5").Value = 2.491514
$ws.Range("N5").Value = 4.983028
$ws.Range("O5").Value = 0.008764135257459093
$ws.Range("P5").Value = 0.005859875752869413
$ws.Range("Q5").Value = 2.05485125636
$ws.Range("R5").Value = 12.32910753816
$ws.Range("S5").Value = 0.00002542563858831824
$ws.Range("T5").Value = 0.0000170498986257946

$ws.Range("G6").Value = 0.8247399999999999
$ws.Range("H6").Value = 2.47422
$ws.Range("I6").Value = 0.00290110066097835
$ws.Range("J6").Value = 0.0029096007056883
$ws.Range("M6").Value = 160.0903473333333
$ws.Range("N6").Value = 480.271042
$ws.Range("O6").Value = 0.5631328812292201
$ws.Range("P6").Value = 0.5647828255873993
$ws.Range("Q6").Value = 132.0329130596933
$ws.Range("R6").Value = 1188.29621753724
$ws.Range("S6").Value = 0.001633705173952733
$ws.Range("T6").Value = 0.001643292507889729

$ws.Range("I7").Value = 0.001056260848766244
$ws.Range("J7").Value = 0.001059355627434439
$ws.Range("M7").Value = 0.8247399999999999
$ws.Range("N7").Value = 2.47422
$ws.Range("O7").Value = 0.00290110066097835
$ws.Range("P7").Value = 0.0029096007056883
$ws.Range("Q7").Value = 0.2476523773733333
$ws.Range("R7").Value = 2.22887139636
$ws.Range("S7").Value = 0.000003064319046521302
$ws.Range("T7").Value = 0.000003082301881158114

$ws.Range("I8").Value = 0.001056260848766244
$ws.Range("J8").Value = 0.001059355627434439
$ws.Range("O8").Value = 0.001056260848766244
$ws.Range("P8").Value = 0.001059355627434439
$ws.Range("S8").Value = 0.000001115686980636385
$ws.Range("T8").Value = 0.000001122234345377013

$ws.Range("I9").Value = 0.001056260848766244
$ws.Range("J9").Value = 0.001059355627434439
$ws.Range("M9").Value = 120.5783256666667
$ws.Range("N9").Value = 361.734977
$ws.Range("O9").Value = 0.4241456220035762
$ws.Range("P9").Value = 0.4253883423266084
$ws.Range("Q9").Value = 36.20717924563623
$ws.Range("R9").Value = 325.864613210726
$ws.Range("S9").Value = 0.0004480084146979837
$ws.Range("T9").Value = 0.0004506375342887

$ws.Range("I10").Value = 0.001056260848766244
$ws.Range("J10").Value = 0.001059355627434439
$ws.Range("M10").Value = 2.491514
$ws.Range("N10").Value = 4.983028
$ws.Range("O10").Value = 0.008764135257459093
$ws.Range("P10").Value = 0.005859875752869413
$ws.Range("Q10").Value = 0.7481501629106667
$ws.Range("R10").Value = 4.488900977464001
$ws.Range("S10").Value = 0.000009257212945745902
$ws.Range("T10").Value = 0.000006207692354868831

$ws.Range("I11").Value = 0.001056260848766244
$ws.Range("J11").Value = 0.001059355627434439
$ws.Range("M11").Value = 160.0903473333333
$ws.Range("N11").Value = 480.271042
$ws.Range("O11").Value = 0.5631328812292201
$ws.Range("P11").Value = 0.5647828255873993
$ws.Range("Q11").Value = 48.07182277035511
$ws.Range("R11").Value = 432.646404933196
$ws.Range("S11").Value = 0.0005948152150953563
$ws.Range("T11").Value = 0.0005983058645643346

$ws.Range("G12").Value = 120.5783256666667
$ws.Range("H12").Value = 361.734977
$ws.Range("I12").Value = 0.4241456220035762
$ws.Range("J12").Value = 0.4253883423266084
$ws.Range("M12").Value = 0.8247399999999999
$ws.Range("N12").Value = 2.47422
$ws.Range("O12").Value = 0.00290110066097835
$ws.Range("P12").Value = 0.0029096007056883
$ws.Range("Q12").Value = 99.44576831032666
$ws.Range("R12").Value = 895.0119147929399
$ws.Range("S12").Value = 0.001230489144345648
$ws.Range("T12").Value = 0.001237710221025076

$ws.Range("G13").Value = 120.5783256666667
$ws.Range("H13").Value = 361.734977
$ws.Range("I13").Value = 0.4241456220035762
$ws.Range("J13").Value = 0.4253883423266084
$ws.Range("O13").Value = 0.001056260848766244
$ws.Range("P13").Value = 0.001059355627434439
$ws.Range("Q13").Value = 36.20717924563623
$ws.Range("R13").Value = 325.864613210726
$ws.Range("S13").Value = 0.0004480084146979837
$ws.Range("T13").Value = 0.0004506375342887

$ws.Range("G14").Value = 120.5783256666667
$ws.Range("H14").Value = 361.734977
$ws.Range("I14").Value = 0.4241456220035762
$ws.Range("J14").Value = 0.4253883423266084
$ws.Range("M14").Value = 120.5783256666667
$ws.Range("N14").Value = 361.734977
$ws.Range("O14").Value = 0.4241456220035762
$ws.Range("P14").Value = 0.4253883423266084
$ws.Range("Q14").Value = 14539.13262057673
$ws.Range("R14").Value = 130852.1935851905
$ws.Range("S14").Value = 0.1798995086648005
$ws.Range("T14").Value = 0.1809552417873798

$ws.Range("G15").Value = 120.5783256666667
$ws.Range("H15").Value = 361.734977
$ws.Range("I15").Value = 0.4241456220035762
$ws.Range("J15").Value = 0.4253883423266084
$ws.Range("M15").Value = 2.491514
$ws.Range("N15").Value = 4.983028
$ws.Range("O15").Value = 0.008764135257459093
$ws.Range("P15").Value = 0.005859875752869413
$ws.Range("Q15").Value = 300.4225864950594
$ws.Range("R15").Value = 1802.535518970356
$ws.Range("S15").Value = 0.00371726960009846
$ws.Range("T15").Value = 0.002492722832753006

$ws.Range("G16").Value = 120.5783256666667
$ws.Range("H16").Value = 361.734977
$ws.Range("I16").Value = 0.4241456220035762
$ws.Range("J16").Value = 0.4253883423266084
$ws.Range("M16").Value = 160.0903473333333
$ws.Range("N16").Value = 480.271042
$ws.Range("O16").Value = 0.5631328812292201
$ws.Range("P16").Value = 0.5647828255873993
$ws.Range("Q16").Value = 19303.42603684845
$ws.Range("R16").Value = 173730.8343316361
$ws.Range("S16").Value = 0.2388503461796336
$ws.Range("T16").Value = 0.2402520299511618

$ws.Range("G17").Value = 2.491514
$ws.Range("H17").Value = 4.983028
$ws.Range("I17").Value = 0.008764135257459093
$ws.Range("J17").Value = 0.005859875752869413
$ws.Range("M17").Value = 0.8247399999999999
$ws.Range("N17").Value = 2.47422
$ws.Range("O17").Value = 0.00290110066097835
$ws.Range("P17").Value = 0.0029096007056883
$ws.Range("Q17").Value = 2.05485125636
$ws.Range("R17").Value = 12.32910753816
$ws.Range("S17").Value = 0.00002542563858831824
$ws.Range("T17").Value = 0.0000170498986257946

$ws.Range("G18").Value = 2.491514
$ws.Range("H18").Value = 4.983028
$ws.Range("I18").Value = 0.008764135257459093
$ws.Range("J18").Value = 0.005859875752869413
$ws.Range("O18").Value = 0.001056260848766244
$ws.Range("P18").Value = 0.001059355627434439
$ws.Range("Q18").Value = 0.7481501629106667
$ws.Range("R18").Value = 4.488900977464001
$ws.Range("S18").Value = 0.000009257212945745902
$ws.Range("T18").Value = 0.000006207692354868831

$ws.Range("G19").Value = 2.491514
$ws.Range("H19").Value = 4.983028
$ws.Range("I19").Value = 0.008764135257459093
$ws.Range("J19").Value = 0.005859875752869413
$ws.Range("M19").Value = 120.5783256666667
$ws.Range("N19").Value = 361.734977
$ws.Range("O19").Value = 0.4241456220035762
$ws.Range("P19").Value = 0.4253883423266084
$ws.Range("Q19").Value = 300.4225864950594
$ws.Range("R19").Value = 1802.535518970356
$ws.Range("S19").Value = 0.00371726960009846
$ws.Range("T19").Value = 0.002492722832753006

$ws.Range("G20").Value = 2.491514
$ws.Range("H20").Value = 4.983028
$ws.Range("I20").Value = 0.008764135257459093
$ws.Range("J20").Value = 0.005859875752869413
$ws.Range("M20").Value = 2.491514
$ws.Range("N20").Value = 4.983028
$ws.Range("O20").Value = 0.008764135257459093
$ws.Range("P20").Value = 0.005859875752869413
$ws.Range("Q20").Value = 6.207642012196
$ws.Range("R20").Value = 24.830568048784
$ws.Range("S20").Value = 0.00007681006681103757
$ws.Range("T20").Value = 0.00003433814383906687

$ws.Range("G21").Value = 2.491514
$ws.Range("H21").Value = 4.983028
$ws.Range("I21").Value = 0.008764135257459093
$ws.Range("J21").Value = 0.005859875752869413
$ws.Range("M21").Value = 160.0903473333333
$ws.Range("N21").Value = 480.271042
$ws.Range("O21").Value = 0.5631328812292201
$ws.Range("P21").Value = 0.5647828255873993
$ws.Range("Q21").Value = 398.8673416458627
$ws.Range("R21").Value = 2393.204049875176
$ws.Range("S21").Value = 0.004935372739015532
$ws.Range("T21").Value = 0.003309557185296676

$ws.Range("G22").Value = 160.0903473333333
$ws.Range("H22").Value = 480.271042
$ws.Range("I22").Value = 0.5631328812292201
$ws.Range("J22").Value = 0.5647828255873993
$ws.Range("M22").Value = 0.8247399999999999
$ws.Range("N22").Value = 2.47422
$ws.Range("O22").Value = 0.00290110066097835
$ws.Range("P22").Value = 0.0029096007056883
$ws.Range("Q22").Value = 132.0329130596933
$ws.Range("R22").Value = 1188.29621753724
$ws.Range("S22").Value = 0.001633705173952733
$ws.Range("T22").Value = 0.001643292507889729

$ws.Range("G23").Value = 160.0903473333333
$ws.Range("H23").Value = 480.271042
$ws.Range("I23").Value = 0.5631328812292201
$ws.Range("J23").Value = 0.5647828255873993
$ws.Range("O23").Value = 0.001056260848766244
$ws.Range("P23").Value = 0.001059355627434439
$ws.Range("Q23").Value = 48.07182277035511
$ws.Range("R23").Value = 432.646404933196
$ws.Range("S23").Value = 0.0005948152150953563
$ws.Range("T23").Value = 0.0005983058645643346

$ws.Range("G24").Value = 160.0903473333333
$ws.Range("H24").Value = 480.271042
$ws.Range("I24").Value = 0.5631328812292201
$ws.Range("J24").Value = 0.5647828255873993
$ws.Range("M24").Value = 120.5783256666667
$ws.Range("N24").Value = 361.734977
$ws.Range("O24").Value = 0.4241456220035762
$ws.Range("P24").Value = 0.4253883423266084
$ws.Range("Q24").Value = 19303.42603684845
$ws.Range("R24").Value = 173730.8343316361
$ws.Range("S24").Value = 0.2388503461796336
$ws.Range("T24").Value = 0.2402520299511618

$ws.Range("G25").Value = 160.0903473333333
$ws.Range("H25").Value = 480.271042
$ws.Range("I25").Value = 0.5631328812292201
$ws.Range("J25").Value = 0.5647828255873993
$ws.Range("M25").Value = 2.491514
$ws.Range("N25").Value = 4.983028
$ws.Range("O25").Value = 0.008764135257459093
$ws.Range("P25").Value = 0.005859875752869413
$ws.Range("Q25").Value = 398.8673416458627
$ws.Range("R25").Value = 2393.204049875176
$ws.Range("S25").Value = 0.004935372739015532
$ws.Range("T25").Value = 0.003309557185296676

$ws.Range("G26").Value = 160.0903473333333
$ws.Range("H26").Value = 480.271042
$ws.Range("I26").Value = 0.5631328812292201
$ws.Range("J26").Value = 0.5647828255873993
$ws.Range("M26").Value = 160.0903473333333
$ws.Range("N26").Value = 480.271042
$ws.Range("O26").Value = 0.5631328812292201
$ws.Range("P26").Value = 0.5647828255873993
$ws.Range("Q26").Value = 25628.91930930731
$ws.Range("R26").Value = 230660.2737837658
$ws.Range("S26").Value = 0.3171186419215229
$ws.Range("T26").Value = 0.3189796400784867
